$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "cooling" values in column C (row 3 becomes 3000, all other
# data rows become 1000 instead of the previous flat 200).
$ws.Range("C2").Value  = 1000
$ws.Range("C3").Value  = 3000
$ws.Range("C4").Value  = 1000
$ws.Range("C5").Value  = 1000
$ws.Range("C6").Value  = 1000
$ws.Range("C7").Value  = 1000
$ws.Range("C8").Value  = 1000
$ws.Range("C9").Value  = 1000
$ws.Range("C10").Value = 1000
$ws.Range("C11").Value = 1000

# Reset the custom row heights back to the sheet default (drops the
# per-row ht="14.25" overrides that were present before).
$ws.Rows("1:11").AutoFit()

# Move the active selection down to C4:C11, matching the refreshed view.
$null = $ws.Range("C4:C11").Select()
